$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-format on cells whose new numeric-looking values must remain literal text
# (otherwise Excel auto-converts strings like "1.00" or "17.50" into numbers, dropping
# trailing zeros / the decimal text representation).
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D17", "D19", "D21", "D22", "D24", "D26", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D42", "D44", "D47", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "42.455.49"
$ws.Range("D3").Value = "2.291.95"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "301.38"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "95.62"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.507"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "34.42"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").Value = "18.93"
$ws.Range("E11").Value = "  +4.06%  "
$ws.Range("D12").Value = "0.0786"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "0.119"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "6.76"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "2.649.73"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "2.309.16"
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "42.405.15"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "12.18"
$ws.Range("E19").Value = "  -5.66%  "
$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "67.72"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  +6.46%  "
$ws.Range("D24").Value = "235.44"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  -3.36%  "
$ws.Range("D28").Value = "2.36"
$ws.Range("E28").Value = "  +14.92%  "
$ws.Range("D29").Value = "165.33"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "9.03"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "31.81"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "4.99"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "17.50"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0699"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "2.32"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "20.21"
$ws.Range("E42").Value = "  +12.32%  "
$ws.Range("D43").Value = "1.964.39"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").Value = "10.46"
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "2.515.49"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.82"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "53.05"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "71.15"
$ws.Range("E51").Value = "  -0.12%  "
